$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
